$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 8.929813666666668
$ws.Cells.Item(2, 8).Value = 26.789441
$ws.Cells.Item(2, 9).Value = 0.3579859341865942
$ws.Cells.Item(2, 10).Value = 0.3579859341865942
$ws.Cells.Item(2, 13).Value = 10.92359866666667
$ws.Cells.Item(2, 14).Value = 32.770796
$ws.Cells.Item(2, 15).Value = 0.2236009040380497
$ws.Cells.Item(2, 16).Value = 0.2236009040380497
$ws.Cells.Item(2, 17).Value = 97.54570066278181
$ws.Cells.Item(2, 18).Value = 877.9113059650363
$ws.Cells.Item(2, 19).Value = 0.08004597851702822
$ws.Cells.Item(2, 20).Value = 0.08004597851702822

# Row 3
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 8.929813666666668
$ws.Cells.Item(3, 8).Value = 26.789441
$ws.Cells.Item(3, 9).Value = 0.3579859341865942
$ws.Cells.Item(3, 10).Value = 0.3579859341865942
$ws.Cells.Item(3, 15).Value = 0.4261214970992155
$ws.Cells.Item(3, 16).Value = 0.4261214970992155
$ws.Cells.Item(3, 17).Value = 185.8951339255018
$ws.Cells.Item(3, 18).Value = 1673.056205329516
$ws.Cells.Item(3, 19).Value = 0.1525455022160528
$ws.Cells.Item(3, 20).Value = 0.1525455022160527

# Row 4
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 8.929813666666668
$ws.Cells.Item(4, 8).Value = 26.789441
$ws.Cells.Item(4, 9).Value = 0.3579859341865942
$ws.Cells.Item(4, 10).Value = 0.3579859341865942
$ws.Cells.Item(4, 13).Value = 13.06524766666667
$ws.Cells.Item(4, 14).Value = 39.195743
$ws.Cells.Item(4, 15).Value = 0.2674394472823625
$ws.Cells.Item(4, 16).Value = 0.2674394472823625
$ws.Cells.Item(4, 17).Value = 116.6702271721848
$ws.Cells.Item(4, 18).Value = 1050.032044549663
$ws.Cells.Item(4, 19).Value = 0.09573956037372296
$ws.Cells.Item(4, 20).Value = 0.09573956037372296

# Row 5
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 8.929813666666668
$ws.Cells.Item(5, 8).Value = 26.789441
$ws.Cells.Item(5, 9).Value = 0.3579859341865942
$ws.Cells.Item(5, 10).Value = 0.3579859341865942
$ws.Cells.Item(5, 13).Value = 4.046901
$ws.Cells.Item(5, 14).Value = 12.140703
$ws.Cells.Item(5, 15).Value = 0.0828381515803724
$ws.Cells.Item(5, 16).Value = 0.0828381515803724
$ws.Cells.Item(5, 17).Value = 36.138071857447
$ws.Cells.Item(5, 18).Value = 325.242646717023
$ws.Cells.Item(5, 19).Value = 0.02965489307979031
$ws.Cells.Item(5, 20).Value = 0.02965489307979031

# Row 6
$ws.Cells.Item(6, 9).Value = 0.03080543264277933
$ws.Cells.Item(6, 10).Value = 0.03080543264277933
$ws.Cells.Item(6, 13).Value = 10.92359866666667
$ws.Cells.Item(6, 14).Value = 32.770796
$ws.Cells.Item(6, 15).Value = 0.2236009040380497
$ws.Cells.Item(6, 16).Value = 0.2236009040380497
$ws.Cells.Item(6, 17).Value = 8.394009999828
$ws.Cells.Item(6, 18).Value = 75.546089998452
$ws.Cells.Item(6, 19).Value = 0.006888122588208706
$ws.Cells.Item(6, 20).Value = 0.006888122588208705

# Row 7
$ws.Cells.Item(7, 9).Value = 0.03080543264277933
$ws.Cells.Item(7, 10).Value = 0.03080543264277933
$ws.Cells.Item(7, 15).Value = 0.4261214970992155
$ws.Cells.Item(7, 16).Value = 0.4261214970992155
$ws.Cells.Item(7, 19).Value = 0.01312685707653017
$ws.Cells.Item(7, 20).Value = 0.01312685707653017

# Row 8
$ws.Cells.Item(8, 9).Value = 0.03080543264277933
$ws.Cells.Item(8, 10).Value = 0.03080543264277933
$ws.Cells.Item(8, 13).Value = 13.06524766666667
$ws.Cells.Item(8, 14).Value = 39.195743
$ws.Cells.Item(8, 15).Value = 0.2674394472823625
$ws.Cells.Item(8, 16).Value = 0.2674394472823625
$ws.Cells.Item(8, 17).Value = 10.039715199249
$ws.Cells.Item(8, 18).Value = 90.35743679324099
$ws.Cells.Item(8, 19).Value = 0.008238587879278954
$ws.Cells.Item(8, 20).Value = 0.008238587879278954

# Row 9
$ws.Cells.Item(9, 9).Value = 0.03080543264277933
$ws.Cells.Item(9, 10).Value = 0.03080543264277933
$ws.Cells.Item(9, 13).Value = 4.046901
$ws.Cells.Item(9, 14).Value = 12.140703
$ws.Cells.Item(9, 15).Value = 0.0828381515803724
$ws.Cells.Item(9, 16).Value = 0.0828381515803724
$ws.Cells.Item(9, 17).Value = 3.109756088529
$ws.Cells.Item(9, 18).Value = 27.987804796761
$ws.Cells.Item(9, 19).Value = 0.002551865098761506
$ws.Cells.Item(9, 20).Value = 0.002551865098761506

# Row 10
$ws.Cells.Item(10, 7).Value = 15.246351
$ws.Cells.Item(10, 8).Value = 45.739053
$ws.Cells.Item(10, 9).Value = 0.6112086331706265
$ws.Cells.Item(10, 10).Value = 0.6112086331706265
$ws.Cells.Item(10, 13).Value = 10.92359866666667
$ws.Cells.Item(10, 14).Value = 32.770796
$ws.Cells.Item(10, 15).Value = 0.2236009040380497
$ws.Cells.Item(10, 16).Value = 0.2236009040380497
$ws.Cells.Item(10, 17).Value = 166.545019455132
$ws.Cells.Item(10, 18).Value = 1498.905175096188
$ws.Cells.Item(10, 19).Value = 0.1366668029328128
$ws.Cells.Item(10, 20).Value = 0.1366668029328127

# Row 11
$ws.Cells.Item(11, 7).Value = 15.246351
$ws.Cells.Item(11, 8).Value = 45.739053
$ws.Cells.Item(11, 9).Value = 0.6112086331706265
$ws.Cells.Item(11, 10).Value = 0.6112086331706265
$ws.Cells.Item(11, 15).Value = 0.4261214970992155
$ws.Cells.Item(11, 16).Value = 0.4261214970992155
$ws.Cells.Item(11, 17).Value = 317.388757124892
$ws.Cells.Item(11, 18).Value = 2856.498814124028
$ws.Cells.Item(11, 19).Value = 0.2604491378066326
$ws.Cells.Item(11, 20).Value = 0.2604491378066326

# Row 12
$ws.Cells.Item(12, 7).Value = 15.246351
$ws.Cells.Item(12, 8).Value = 45.739053
$ws.Cells.Item(12, 9).Value = 0.6112086331706265
$ws.Cells.Item(12, 10).Value = 0.6112086331706265
$ws.Cells.Item(12, 13).Value = 13.06524766666667
$ws.Cells.Item(12, 14).Value = 39.195743
$ws.Cells.Item(12, 15).Value = 0.2674394472823625
$ws.Cells.Item(12, 16).Value = 0.2674394472823625
$ws.Cells.Item(12, 17).Value = 199.197351827931
$ws.Cells.Item(12, 18).Value = 1792.776166451379
$ws.Cells.Item(12, 19).Value = 0.1634612990293606
$ws.Cells.Item(12, 20).Value = 0.1634612990293606

# Row 13
$ws.Cells.Item(13, 7).Value = 15.246351
$ws.Cells.Item(13, 8).Value = 45.739053
$ws.Cells.Item(13, 9).Value = 0.6112086331706265
$ws.Cells.Item(13, 10).Value = 0.6112086331706265
$ws.Cells.Item(13, 13).Value = 4.046901
$ws.Cells.Item(13, 14).Value = 12.140703
$ws.Cells.Item(13, 15).Value = 0.0828381515803724
$ws.Cells.Item(13, 16).Value = 0.0828381515803724
$ws.Cells.Item(13, 17).Value = 61.700473108251
$ws.Cells.Item(13, 18).Value = 555.304257974259
$ws.Cells.Item(13, 19).Value = 0.05063139340182059
$ws.Cells.Item(13, 20).Value = 0.05063139340182059
